$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 02:45:26"
$wsZh.Range("H2").Value = "2016-03-12 02:45:38"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 02:45:29"
$wsDe.Range("H2").Value = "2016-03-12 02:45:43"
